$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'308.56"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'0.15%"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'40.81"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'1.72%"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.113"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'0.11%"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'0.07636"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'-1.49%"
$c.Style = "Normal"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D6")
$c.Value = "'1.609"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'-1.12%"
$c.Style = "Normal"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c = $ws.Range("D7")
$c.Value = "'2.452"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'0.68%"
$c.Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D8")
$c.Value = "'0.9043"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'2.65%"
$c.Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c = $ws.Range("D9")
$c.Value = "'0.1110"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'8.34%"
$c.Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Range("D10")
$c.Value = "'0.1782"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'2.11%"
$c.Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Range("D11")
$c.Value = "'0.09233"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'3.09%"
$c.Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Range("D12")
$c.Value = "'0.04198"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'-5.17%"
$c.Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Range("D13")
$c.Value = "'0.1054"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'-0.18%"
$c.Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Range("D14")
$c.Value = "'0.001258"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'-0.51%"
$c.Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c = $ws.Range("D15")
$c.Value = "'0.005684"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'-2.04%"
$c.Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D16")
$c.Value = "'3.351"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'-0.08%"
$c.Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c = $ws.Range("D17")
$c.Value = "'4.246"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'-0.18%"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'0.34%"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'6.561"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'-6.37%"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'1.81%"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'0.2831"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'1.59%"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'0.04067"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'-2.70%"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'2.40%"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'0.004111"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'0.63%"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'0.0001302"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'0.10%"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'0.0003749"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'25.74%"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'0.02422"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'2.40%"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.05193"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'-0.35%"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.007794"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'-1.80%"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.1302"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'-1.71%"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.007050"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'11.16%"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'-0.61%"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.008792"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'0.33%"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.3328"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'-0.51%"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.00006939"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'6.14%"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'0.09%"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'0.03120"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'1,042.71%"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.00002103"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'0.09%"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'0.09%"
$c.Style = "Normal"
